$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 99950
$ws.Range("E2").Value = 4270
$ws.Range("F2").Value = 4270
$ws.Range("G2").Value = 1838
$ws.Range("H2").Value = 1297
$ws.Range("I2").Value = 1324
$ws.Range("J2").Value = -27
$ws.Range("K2").Value = 102602
$ws.Range("L2").Value = 75346
$ws.Range("M2").Value = 27256
$ws.Range("N2").Value = 26297
$ws.Range("O2").Value = 959
$ws.Range("P2").Value = 20781
$ws.Range("Q2").Value = 2644
$ws.Range("R2").Value = -607
$ws.Range("S2").Value = -2891
$ws.Range("T2").Value = 1560
$ws.Range("U2").Value = 1084
$ws.Range("V2").Value = 26857
$ws.Range("W2").Value = 4.27
$ws.Range("X2").Value = 1.3
$ws.Range("Y2").Value = 5.07
$ws.Range("Z2").Value = 1.27
$ws.Range("AA2").Value = 276.44
$ws.Range("AB2").Value = 38.77
$ws.Range("AC2").Value = 319
$ws.Range("AD2").Value = 18.49
$ws.Range("AE2").Value = 6400
$ws.Range("AF2").Value = 0.92
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 415622638

# Row 3
$ws.Range("D3").Value = 98900
$ws.Range("E3").Value = 1689
$ws.Range("F3").Value = 3434
$ws.Range("G3").Value = 1195
$ws.Range("H3").Value = 1046
$ws.Range("I3").Value = 1059
$ws.Range("J3").Value = -13
$ws.Range("K3").Value = 100637
$ws.Range("L3").Value = 72704
$ws.Range("M3").Value = 27933
$ws.Range("N3").Value = 27355
$ws.Range("O3").Value = 578
$ws.Range("P3").Value = 20781
$ws.Range("Q3").Value = 6850
$ws.Range("R3").Value = -5411
$ws.Range("S3").Value = 690
$ws.Range("T3").Value = 4406
$ws.Range("U3").Value = 2445
$ws.Range("V3").Value = 25340
$ws.Range("W3").Value = 1.71
$ws.Range("X3").Value = 1.06
$ws.Range("Y3").Value = 3.95
$ws.Range("Z3").Value = 1.03
$ws.Range("AA3").Value = 260.28
$ws.Range("AB3").Value = 43.67
$ws.Range("AC3").Value = 255
$ws.Range("AD3").Value = 21.98
$ws.Range("AE3").Value = 6658
$ws.Range("AF3").Value = 0.84
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 415622638

# Row 4
$ws.Range("D4").Value = 111059
$ws.Range("E4").Value = -4672
$ws.Range("F4").Value = -4672
$ws.Range("G4").Value = -10111
$ws.Range("H4").Value = -7549
$ws.Range("I4").Value = -7358
$ws.Range("J4").Value = -191
$ws.Range("K4").Value = 99702
$ws.Range("L4").Value = 79003
$ws.Range("M4").Value = 20699
$ws.Range("N4").Value = 20316
$ws.Range("O4").Value = 383
$ws.Range("P4").Value = 20781
$ws.Range("Q4").Value = 3588
$ws.Range("R4").Value = -2618
$ws.Range("S4").Value = 1745
$ws.Range("T4").Value = 2736
$ws.Range("U4").Value = 852
$ws.Range("V4").Value = 27499
$ws.Range("W4").Value = -4.21
$ws.Range("X4").Value = -6.8
$ws.Range("Y4").Value = -30.87
$ws.Range("Z4").Value = -7.54
$ws.Range("AA4").Value = 381.68
$ws.Range("AB4").Value = 8.98
$ws.Range("AC4").Value = -1770
$ws.Range("AD4").Value = -2.87
$ws.Range("AE4").Value = 4944
$ws.Range("AF4").ClearContents() | Out-Null
$ws.Range("AG4").ClearContents() | Out-Null
$ws.Range("AH4").Value = 1.03
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 415622638

# Row 5
$ws.Range("D5").Value = 117668
$ws.Range("E5").Value = 4290
$ws.Range("F5").Value = 4290
$ws.Range("G5").Value = 3136
$ws.Range("H5").Value = 2579
$ws.Range("I5").Value = 2589
$ws.Range("J5").Value = -10
$ws.Range("K5").Value = 87763
$ws.Range("L5").Value = 64983
$ws.Range("M5").Value = 22781
$ws.Range("N5").Value = 22426
$ws.Range("O5").Value = 355
$ws.Range("P5").Value = 20781
$ws.Range("Q5").Value = 3374
$ws.Range("R5").Value = 849
$ws.Range("S5").Value = -7211
$ws.Range("T5").Value = 680
$ws.Range("U5").Value = 2694
$ws.Range("V5").Value = 19705
$ws.Range("W5").Value = 3.65
$ws.Range("X5").Value = 2.19
$ws.Range("Y5").Value = 12.12
$ws.Range("Z5").Value = 2.75
$ws.Range("AA5").Value = 285.25
$ws.Range("AB5").Value = 22.01
$ws.Range("AC5").Value = 623
$ws.Range("AD5").Value = 9.5
$ws.Range("AE5").Value = 5458
$ws.Range("AF5").ClearContents() | Out-Null
$ws.Range("AG5").ClearContents() | Out-Null
$ws.Range("AH5").Value = 1.08
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 415622638

# Row 6
$ws.Range("D6").Value = 106055
$ws.Range("E6").Value = 6287
$ws.Range("F6").Value = 6287
$ws.Range("G6").Value = 4318
$ws.Range("H6").Value = 2973
$ws.Range("I6").Value = 2987
$ws.Range("K6").Value = 87338
$ws.Range("L6").Value = 64161
$ws.Range("M6").Value = 23177
$ws.Range("N6").Value = 22819
$ws.Range("P6").Value = 20781
$ws.Range("Q6").Value = 1762
$ws.Range("R6").Value = -1528
$ws.Range("S6").Value = 1675
$ws.Range("T6").Value = 270
$ws.Range("U6").Value = 1492
$ws.Range("V6").Value = 21470
$ws.Range("W6").Value = 5.93
$ws.Range("X6").Value = 2.8
$ws.Range("Y6").Value = 13.2
$ws.Range("Z6").Value = 3.4
$ws.Range("AA6").Value = 276.83
$ws.Range("AB6").Value = 25.78
$ws.Range("AC6").Value = 719
$ws.Range("AD6").Value = 7.5
$ws.Range("AE6").Value = 5554
$ws.Range("AF6").ClearContents() | Out-Null
$ws.Range("AG6").ClearContents() | Out-Null
$ws.Range("AH6").ClearContents() | Out-Null
$ws.Range("AI6").Value = 0.97
$ws.Range("AJ6").Value = 415622638

# Row 7
$ws.Range("D7").Value = 85901
$ws.Range("E7").Value = 4342
$ws.Range("G7").Value = 3433
$ws.Range("H7").Value = 2383
$ws.Range("I7").Value = 2438
$ws.Range("K7").Value = 91939
$ws.Range("L7").Value = 66217
$ws.Range("M7").Value = 25723
$ws.Range("N7").Value = 25414
$ws.Range("P7").Value = 20780
$ws.Range("Q7").Value = 2135
$ws.Range("R7").Value = -1111
$ws.Range("S7").Value = 1507
$ws.Range("T7").Value = 395
$ws.Range("U7").Value = 1449
$ws.Range("W7").Value = 5.05
$ws.Range("X7").Value = 2.77
$ws.Range("Y7").Value = 10.11
$ws.Range("Z7").Value = 2.66
$ws.Range("AA7").Value = 257.42
$ws.Range("AC7").Value = 587
$ws.Range("AD7").Value = 7.96
$ws.Range("AE7").Value = 6185
$ws.Range("AF7").Value = 0.76
$ws.Range("AG7").Value = 0
$ws.Range("AH7").ClearContents() | Out-Null
$ws.Range("AI7").Value = 0

# Row 8
$ws.Range("D8").Value = 92158
$ws.Range("E8").Value = 5306
$ws.Range("G8").Value = 4293
$ws.Range("H8").Value = 3162
$ws.Range("I8").Value = 3177
$ws.Range("K8").Value = 95977
$ws.Range("L8").Value = 67101
$ws.Range("M8").Value = 28875
$ws.Range("N8").Value = 28620
$ws.Range("P8").Value = 20780
$ws.Range("Q8").Value = 3732
$ws.Range("R8").Value = -857
$ws.Range("S8").Value = -971
$ws.Range("T8").Value = 428
$ws.Range("U8").Value = 2952
$ws.Range("W8").Value = 5.76
$ws.Range("X8").Value = 3.43
$ws.Range("Y8").Value = 11.76
$ws.Range("Z8").Value = 3.37
$ws.Range("AA8").Value = 232.39
$ws.Range("AC8").Value = 764
$ws.Range("AD8").Value = 5.69
$ws.Range("AE8").Value = 6965
$ws.Range("AF8").Value = 0.62
$ws.Range("AG8").Value = 0
$ws.Range("AH8").ClearContents() | Out-Null
$ws.Range("AI8").Value = 0

# Row 9
$ws.Range("D9").Value = 99807
$ws.Range("E9").Value = 5853
$ws.Range("G9").Value = 4893
$ws.Range("H9").Value = 3605
$ws.Range("I9").Value = 3670
$ws.Range("K9").Value = 100658
$ws.Range("L9").Value = 68186
$ws.Range("M9").Value = 32473
$ws.Range("N9").Value = 32266
$ws.Range("P9").Value = 20780
$ws.Range("Q9").Value = 2855
$ws.Range("R9").Value = -1652
$ws.Range("S9").Value = 94
$ws.Range("T9").Value = 408
$ws.Range("U9").Value = 2509
$ws.Range("W9").Value = 5.86
$ws.Range("X9").Value = 3.61
$ws.Range("Y9").Value = 12.06
$ws.Range("Z9").Value = 3.67
$ws.Range("AA9").Value = 209.98
$ws.Range("AC9").Value = 883
$ws.Range("AD9").Value = 4.93
$ws.Range("AE9").Value = 7853
$ws.Range("AF9").Value = 0.55
$ws.Range("AG9").Value = 6
$ws.Range("AH9").Value = 0.14
$ws.Range("AI9").Value = 0.67
